$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C11").Value = 'Los Lagos'
$ws.Range("D11").Value = 44950
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112031
$ws.Range("G11").Value = 'Poroto verde'
$ws.Range("H11").Value = 'Magnum'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 40
$ws.Range("K11").Value = 25000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 25000
$ws.Range("N11").Value = '$/saco 25 kilos'
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("P11").Value = 1000
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = 'Hortaliza'
